$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row for "Make your own Neural Network" by Tariq Rashid
$ws.Range("A17").Value = "Make your own Neural Network / Neuronale Netze selbst programmieren"
$ws.Range("B17").Value = "Tariq Rashid"
$ws.Range("C17").Value = 211
$ws.Range("D17").Value = "x"

# Update the selection to reflect where the user ended up
$ws.Range("E14").Select()
